$wb = $excel.ActiveWorkbook

foreach ($name in @("two_blank_rows", "occupied_row_and_blank_row", "two_occupied_rows")) {
    $ws = $wb.Worksheets.Item($name)

    # Clear out the old row 4 (becomes an embedded blank row)
    $ws.Range("A4:B4").ClearContents()

    # Replace row 5's numeric values with string values
    $ws.Range("A5").Value = "v2,1"
    $ws.Range("B5").Value = "v2,2"

    # Add a brand new row 7 (row 6 stays blank)
    $ws.Range("A7").Value = "v4,1"
    $ws.Range("B7").Value = "v4,2"

    # Update the visible selection to match
    $ws.Range("A3:B7").Select()
}

# Restore the originally active sheet/tab
$wb.Worksheets.Item("two_occupied_rows").Activate()
